$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B to hold the full "Skill Description"
# (the long/full name of the skill), shifting the old B/C/D (SFIA Level,
# Keycode, Description) one column to the right (C/D/E).
$ws.Columns.Item(2).Insert()

$ws.Range("B1").Value = "Skill Description"

# Rows 2-3: Autonomy
$ws.Range("B2").Value = "Autonomy"
$ws.Range("B3").Value = "Autonomy"

# Rows 4-7: Influence
$ws.Range("B4").Value = "Influence"
$ws.Range("B5").Value = "Influence"
$ws.Range("B6").Value = "Influence"
$ws.Range("B7").Value = "Influence"

# Rows 8-10: Complexity
$ws.Range("B8").Value = "Complexity"
$ws.Range("B9").Value = "Complexity"
$ws.Range("B10").Value = "Complexity"

# Rows 11-13: Knowledge
$ws.Range("B11").Value = "Knowledge"
$ws.Range("B12").Value = "Knowledge"
$ws.Range("B13").Value = "Knowledge"

# Rows 15-18: BURM -> Risk management
$ws.Range("B15").Value = "Risk management"
$ws.Range("B16").Value = "Risk management"
$ws.Range("B17").Value = "Risk management"
$ws.Range("B18").Value = "Risk management"
